$d = $word.ActiveDocument

# Update the date/day heading
$d.Content.Find.Execute("2023-12-14 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-15 Friday", 1)

# Update the division problems in the table, cell by cell (so the two
# "69÷9=" occurrences each get their own distinct replacement).
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$replacements = @(
    @("42÷9=", "97÷6="), @("26÷5=", "69÷3="), @("70÷3=", "87÷5="), @("14÷5=", "25÷9="), @("78÷9=", "38÷8="),
    @("74÷9=", "52÷2="), @("46÷4=", "15÷2="), @("60÷6=", "40÷7="), @("96÷4=", "77÷7="), @("21÷6=", "35÷4="),
    @("70÷9=", "73÷4="), @("37÷3=", "82÷8="), @("15÷7=", "41÷4="), @("68÷5=", "20÷2="), @("76÷5=", "69÷5="),
    @("45÷8=", "38÷9="), @("67÷5=", "95÷3="), @("75÷8=", "17÷2="), @("53÷7=", "92÷4="), @("69÷9=", "62÷2="),
    @("72÷8=", "58÷3="), @("12÷3=", "35÷6="), @("87÷2=", "43÷6="), @("31÷8=", "56÷7="), @("69÷9=", "99÷6=")
)

$idx = 0
foreach ($r in $rows) {
    for ($c = 1; $c -le 5; $c++) {
        $old = $replacements[$idx][0]
        $new = $replacements[$idx][1]
        $cellRange = $t.Cell($r, $c).Range
        $cellRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
        $idx = $idx + 1
    }
}
